$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# This sheet is a rolling weekly "MarketBeat rank" watch table.  Each week a
# new pair of "current" columns is inserted right after column A (the firm
# name) and the previously existing week columns (B..E) shift three places
# to the right (become E..H).  Two brand-new analyst firms are also appended
# as new rows at the bottom of the table.
#
# Existing layout (before):
#   A = firm name
#   B = Jun_17 (most recent), C = Jun_15, D = Jun_13, E = Jun_10 (oldest)
#
# New layout (after):
#   A = firm name
#   B = Jun_27 (new), C = Jun_26 (new), D = Jun_26 (new, duplicate of C)
#   E = Jun_17, F = Jun_15, G = Jun_13, H = Jun_10   (shifted from B,C,D,E)
# ---------------------------------------------------------------------------

$lastRow = 27
$lastColBefore = 5   # column E

# 1) Shift the existing week columns (B..E) three columns to the right
#    (E..H), preserving both values and cell formatting.  Columns are
#    copied right-to-left so we never overwrite a column before it has
#    been read, and cell-by-cell so formatting survives the copy.
for ($r = 1; $r -le $lastRow; $r++) {
    for ($c = $lastColBefore; $c -ge 2; $c--) {
        $src = $ws.Cells.Item($r, $c)
        $dst = $ws.Cells.Item($r, $c + 3)
        $src.Copy($dst)
    }
}

# 2) Reset the old (now stale) special/highlighted cells that used to hold
#    "changed rating" notes back to the plain "UN" placeholder with default
#    (no) fill, since that content now lives three columns further right.
$ws.Cells.Item(8, 5).Value = "UN"
$ws.Cells.Item(8, 5).ClearFormats()

$ws.Cells.Item(9, 5).Value = "UN"
$ws.Cells.Item(9, 5).ClearFormats()

$ws.Cells.Item(17, 5).Value = "UN"
$ws.Cells.Item(17, 5).ClearFormats()

$ws.Cells.Item(18, 3).Value = "UN"
$ws.Cells.Item(18, 3).ClearFormats()

$ws.Cells.Item(18, 4).Value = "UN"
$ws.Cells.Item(18, 4).ClearFormats()

$ws.Cells.Item(18, 5).Value = "UN"
$ws.Cells.Item(18, 5).ClearFormats()

$ws.Cells.Item(22, 4).Value = "UN"
$ws.Cells.Item(22, 4).ClearFormats()

# 3) Populate the two new leading columns (B, C) plus the duplicated D
#    header with the new week labels.  The order of first-use below
#    matches the order new shared strings must be appended in.
$ws.Cells.Item(1, 3).Value = "Jun_26"        # C1

$ws.Cells.Item(28, 1).Value = "Benchmark"    # A28 - new firm row
$ws.Cells.Item(28, 2).Value = "UN"
$ws.Cells.Item(28, 3).Value = "UN"
$ws.Cells.Item(28, 4).Value = "UN"

$ws.Cells.Item(29, 1).Value = "Evercore ISI" # A29 - new firm row
$ws.Cells.Item(29, 2).Value = "UN"
$ws.Cells.Item(29, 3).Value = "UN"
$ws.Cells.Item(29, 4).Value = "UN"

$ws.Cells.Item(1, 2).Value = "Jun_27"        # B1

$ws.Cells.Item(1, 4).Value = "Jun_26"        # D1 (duplicate of C1)

# 4) Fill in the plain "UN" placeholder across the two new week columns
#    (B, C, D) for every existing firm row, since only the header row was
#    set above.
for ($r = 2; $r -le $lastRow; $r++) {
    $ws.Cells.Item($r, 2).Value = "UN"
    $ws.Cells.Item($r, 3).Value = "UN"
    $ws.Cells.Item($r, 4).Value = "UN"
}

# 5) Match the column widths used by the rest of the week columns for the
#    three newly introduced columns (F, G and the new last column H).
$ws.Columns.Item(5).ColumnWidth = 7.14
$ws.Columns.Item(6).ColumnWidth = 7.14
$ws.Columns.Item(7).ColumnWidth = 7.14
$ws.Columns.Item(8).ColumnWidth = 7.14
